$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9220307469367981
$ws.Range("B1").Value = 1.432577133178711
$ws.Range("C1").Value = 4.889541625976562
$ws.Range("D1").Value = 2.871509552001953
$ws.Range("E1").Value = 0.4070788621902466
